# Auto-generated edit script: applies scheduled-runner price/profit updates
# to the Asura_Profits workbook (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 289.66666
$ws.Range("I2").Value = 297.6
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 297.6
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -184.6
$ws.Range("N2").Value = -476
$ws.Range("H17").Value = 25247.414
$ws.Range("J17").Value = 25247.414
$ws.Range("L17").Value = 75742.242
$ws.Range("N17").Value = -76078.242
$ws.Range("H40").Value = 1581.1666
$ws.Range("I40").Value = 1436.1
$ws.Range("J40").Value = 1762.5
$ws.Range("K40").Value = 1436.1
$ws.Range("L40").Value = 1762.5
$ws.Range("M40").Value = -1261.1
$ws.Range("N40").Value = -2112.5
$ws.Range("H58").Value = 2183.0588
$ws.Range("J58").Value = 2478.3572
$ws.Range("L58").Value = 7435.071599999999
$ws.Range("N58").Value = -7735.071599999999
$ws.Range("H86").Value = 2522
$ws.Range("I86").Value = 3083.3333
$ws.Range("J86").Value = 2101
$ws.Range("K86").Value = 3083.3333
$ws.Range("L86").Value = 2101
$ws.Range("M86").Value = -1960.3333
$ws.Range("N86").Value = -4347
$ws.Range("H89").Value = 2522
$ws.Range("I89").Value = 3083.3333
$ws.Range("J89").Value = 2101
$ws.Range("K89").Value = 15416.6665
$ws.Range("L89").Value = 10505
$ws.Range("M89").Value = -9800.666499999999
$ws.Range("N89").Value = -21737
$ws.Range("H92").Value = 7682.0713
$ws.Range("I92").Value = 16857.334
$ws.Range("J92").Value = 800.625
$ws.Range("K92").Value = 16857.334
$ws.Range("L92").Value = 800.625
$ws.Range("M92").Value = -15609.334
$ws.Range("N92").Value = -3296.625
$ws.Range("H99").Value = 543.3333
$ws.Range("I99").Value = 553
$ws.Range("J99").Value = 495
$ws.Range("K99").Value = 1659
$ws.Range("L99").Value = 1485
$ws.Range("M99").Value = -161
$ws.Range("N99").Value = -4481
$ws.Range("H100").Value = 2436.8147
$ws.Range("I100").Value = 1649.9166
$ws.Range("J100").Value = 3066.3333
$ws.Range("K100").Value = 1649.9166
$ws.Range("L100").Value = 3066.3333
$ws.Range("M100").Value = -1108.9166
$ws.Range("N100").Value = -4148.3333
$ws.Range("H138").Value = 2932.2917
$ws.Range("I138").Value = 1958.25
$ws.Range("J138").Value = 4462.9287
$ws.Range("K138").Value = 5874.75
$ws.Range("L138").Value = 13388.7861
$ws.Range("M138").Value = -734.75
$ws.Range("N138").Value = -23668.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1153
$ws.Range("I61").Value = 996.8929000000001
$ws.Range("J61").Value = 2245.75
$ws.Range("K61").Value = 996.8929000000001
$ws.Range("L61").Value = 2245.75
$ws.Range("M61").Value = -784.8929000000001
$ws.Range("N61").Value = -2669.75
$ws.Range("H136").Value = 1153
$ws.Range("I136").Value = 996.8929000000001
$ws.Range("J136").Value = 2245.75
$ws.Range("K136").Value = 2990.6787
$ws.Range("L136").Value = 6737.25
$ws.Range("M136").Value = -440.6787000000004
$ws.Range("N136").Value = -11837.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 336333.34
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 336333.34
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -337235.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1825.8959
$ws.Range("I31").Value = 1295.8108
$ws.Range("J31").Value = 3608.9092
$ws.Range("K31").Value = 1295.8108
$ws.Range("L31").Value = 3608.9092
$ws.Range("M31").Value = -1000.8108
$ws.Range("N31").Value = -4198.9092
$ws.Range("H34").Value = 1825.8959
$ws.Range("I34").Value = 1295.8108
$ws.Range("J34").Value = 3608.9092
$ws.Range("K34").Value = 1295.8108
$ws.Range("L34").Value = 3608.9092
$ws.Range("M34").Value = -1093.8108
$ws.Range("N34").Value = -4012.9092
$ws.Range("H62").Value = 85800.664
$ws.Range("I62").Value = 168835
$ws.Range("J62").Value = 2766.3333
$ws.Range("K62").Value = 168835
$ws.Range("L62").Value = 2766.3333
$ws.Range("M62").Value = -168211
$ws.Range("N62").Value = -4014.3333
$ws.Range("H65").Value = 85800.664
$ws.Range("I65").Value = 168835
$ws.Range("J65").Value = 2766.3333
$ws.Range("K65").Value = 844175
$ws.Range("L65").Value = 13831.6665
$ws.Range("M65").Value = -841055
$ws.Range("N65").Value = -20071.6665
$ws.Range("H105").Value = 8397.071
$ws.Range("I105").Value = 8966.076999999999
$ws.Range("K105").Value = 8966.076999999999
$ws.Range("M105").Value = -7219.076999999999
$ws.Range("H132").Value = 502445.16
$ws.Range("I132").Value = 588804.6
$ws.Range("J132").Value = 5878
$ws.Range("K132").Value = 1766413.8
$ws.Range("L132").Value = 17634
$ws.Range("M132").Value = -1763883.8
$ws.Range("N132").Value = -22694

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22751910
$ws.Range("I131").Value = 11273
$ws.Range("J131").Value = 29440334
$ws.Range("K131").Value = 33819
$ws.Range("L131").Value = 88321002
$ws.Range("M131").Value = -28779
$ws.Range("N131").Value = -88331082

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2188.1667
$ws.Range("I126").Value = 1952
$ws.Range("J126").Value = 2518.8
$ws.Range("K126").Value = 5856
$ws.Range("L126").Value = 7556.400000000001
$ws.Range("M126").Value = -3386
$ws.Range("N126").Value = -12496.4
$ws.Range("H132").Value = 1766.3334
$ws.Range("I132").Value = 922.5294
$ws.Range("J132").Value = 3815.5715
$ws.Range("K132").Value = 2767.5882
$ws.Range("L132").Value = 11446.7145
$ws.Range("M132").Value = -237.5882000000001
$ws.Range("N132").Value = -16506.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7078.25
$ws.Range("I23").Value = 5153
$ws.Range("K23").Value = 5153
$ws.Range("M23").Value = -4923
$ws.Range("H40").Value = 2793.182
$ws.Range("I40").Value = 2215.625
$ws.Range("K40").Value = 2215.625
$ws.Range("M40").Value = -2079.625
$ws.Range("H68").Value = 2056.3076
$ws.Range("I68").Value = 1434.5333
$ws.Range("J68").Value = 2904.182
$ws.Range("K68").Value = 1434.5333
$ws.Range("L68").Value = 2904.182
$ws.Range("M68").Value = -685.5333000000001
$ws.Range("N68").Value = -4402.182
$ws.Range("H71").Value = 2056.3076
$ws.Range("I71").Value = 1434.5333
$ws.Range("J71").Value = 2904.182
$ws.Range("K71").Value = 7172.6665
$ws.Range("L71").Value = 14520.91
$ws.Range("M71").Value = -3428.6665
$ws.Range("N71").Value = -22008.91
$ws.Range("H93").Value = 1275.375
$ws.Range("I93").Value = 840.6
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 840.6
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 407.4
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 9013.25
$ws.Range("I61").Value = 5526
$ws.Range("J61").Value = 12500.5
$ws.Range("K61").Value = 5526
$ws.Range("L61").Value = 12500.5
$ws.Range("M61").Value = -5234
$ws.Range("N61").Value = -13084.5
$ws.Range("H96").Value = 2501.5
$ws.Range("J96").Value = 2500
$ws.Range("L96").Value = 2500
$ws.Range("N96").Value = -5246
$ws.Range("H132").Value = 1674.4
$ws.Range("I132").Value = 1317.5555
$ws.Range("J132").Value = 2209.6667
$ws.Range("K132").Value = 3952.6665
$ws.Range("L132").Value = 6629.000100000001
$ws.Range("M132").Value = -1422.6665
$ws.Range("N132").Value = -11689.0001
$ws.Range("H136").Value = 1273.619
$ws.Range("I136").Value = 986.6316
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2959.8948
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -409.8948
$ws.Range("N136").Value = -17100
